$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.735.56'
$ws.Range('E2').Value = '  -0.16%  '

$ws.Range('D3').Value = '3.235.26'
$ws.Range('E3').Value = '  -0.58%  '

$ws.Range('E4').Value = '  +0.02%  '

$c = $ws.Range('D5')
$c.NumberFormat = "@"
$c.Value = '579.36'
$c.ClearFormats()
$ws.Range('E5').Value = '  -0.78%  '

$c = $ws.Range('D6')
$c.NumberFormat = "@"
$c.Value = '183.48'
$c.ClearFormats()
$ws.Range('E6').Value = '  +0.04%  '

$ws.Range('E7').Value = '  +0.02%  '

$ws.Range('E8').Value = '  +0.31%  '

$c = $ws.Range('D9')
$c.NumberFormat = "@"
$c.Value = '0.130'
$c.ClearFormats()
$ws.Range('E9').Value = '  -3.67%  '

$ws.Range('E10').Value = '  -1.08%  '

$c = $ws.Range('D11')
$c.NumberFormat = "@"
$c.Value = '0.414'
$c.ClearFormats()
$ws.Range('E11').Value = '  -0.31%  '

$ws.Range('D12').Value = '3.799.35'
$ws.Range('E12').Value = '  -0.57%  '

$ws.Range('E13').Value = '  +0.00%  '

$ws.Range('E14').Value = '  -3.58%  '

$ws.Range('D15').Value = '67.763.48'
$ws.Range('E15').Value = '  -0.12%  '

$ws.Range('E16').Value = '  -1.72%  '

$ws.Range('D17').Value = '3.268.75'
$ws.Range('E17').Value = '  +0.68%  '

$ws.Range('E18').Value = '  -1.03%  '

$c = $ws.Range('D19')
$c.NumberFormat = "@"
$c.Value = '13.43'
$c.ClearFormats()
$ws.Range('E19').Value = '  -1.04%  '

$c = $ws.Range('D20')
$c.NumberFormat = "@"
$c.Value = '395.43'
$c.ClearFormats()
$ws.Range('E20').Value = '  +3.68%  '

$c = $ws.Range('D21')
$c.NumberFormat = "@"
$c.Value = '7.55'
$c.ClearFormats()
$ws.Range('E21').Value = '  -1.30%  '

$ws.Range('E22').Value = '  +0.02%  '

$c = $ws.Range('D23')
$c.NumberFormat = "@"
$c.Value = '71.07'
$c.ClearFormats()

$ws.Range('E24').Value = '  -0.04%  '

$ws.Range('E25').Value = '  -1.68%  '

$ws.Range('E26').Value = '  +2.24%  '

$ws.Range('E27').Value = '  -3.10%  '

$ws.Range('E28').Value = '  -0.22%  '

$ws.Range('E29').Value = '  -1.32%  '

$ws.Range('E30').Value = '  -2.02%  '

$c = $ws.Range('D31')
$c.NumberFormat = "@"
$c.Value = '22.62'
$c.ClearFormats()
$ws.Range('E31').Value = '  -1.24%  '

$ws.Range('E32').Value = '  -2.52%  '

$ws.Range('E33').Value = '  -1.26%  '

$c = $ws.Range('D34')
$c.NumberFormat = "@"
$c.Value = '0.999'
$c.ClearFormats()
$ws.Range('E34').Value = '  +0.02%  '

$c = $ws.Range('D35')
$c.NumberFormat = "@"
$c.Value = '161.95'
$c.ClearFormats()
$ws.Range('E35').Value = '  +0.38%  '

$ws.Range('E36').Value = '  -4.77%  '

$ws.Range('E37').Value = '  +1.50%  '

$c = $ws.Range('D39')
$c.NumberFormat = "@"
$c.Value = '26.35'
$c.ClearFormats()
$ws.Range('E39').Value = '  -0.88%  '

$ws.Range('E40').Value = '  -1.19%  '

$ws.Range('E41').Value = '  -3.28%  '

$c = $ws.Range('D42')
$c.NumberFormat = "@"
$c.Value = '41.11'
$c.ClearFormats()
$ws.Range('E42').Value = '  -0.62%  '

$ws.Range('E43').Value = '  -4.98%  '

$c = $ws.Range('D44')
$c.NumberFormat = "@"
$c.Value = '0.0683'
$c.ClearFormats()
$ws.Range('E44').Value = '  -0.69%  '

$ws.Range('D45').Value = '2.604.04'

$c = $ws.Range('D46')
$c.NumberFormat = "@"
$c.Value = '24.84'
$c.ClearFormats()
$ws.Range('E46').Value = '  -2.33%  '

$c = $ws.Range('D47')
$c.NumberFormat = "@"
$c.Value = '335.21'
$c.ClearFormats()
$ws.Range('E47').Value = '  -3.48%  '

$ws.Range('E48').Value = '  -2.34%  '

$ws.Range('E49').Value = '  +0.99%  '

$ws.Range('E50').Value = '  -2.07%  '

$ws.Range('B51').Value = 'ONDO'
$ws.Range('C51').Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$c = $ws.Range('D51')
$c.NumberFormat = "@"
$c.Value = '0.971'
$c.ClearFormats()
$ws.Range('E51').Value = '  -2.09%  '

Write-Output "applied"
